# Regenerate save_data to use K (strikeouts) instead of Strike# column values.
# Update the "K" column (column G) values for the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 0
    6  = 1
    7  = 1
    8  = 2
    9  = 2
    10 = 2
    11 = 2
    12 = 1
    13 = 3
    14 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
